$wb = $excel.ActiveWorkbook

# New headers row (common to every sheet): eb, gb, hp, st, wi, ieh, chp, ac, ab_ct, ab_hp, cp_ct, cp_hp, ttes, btes, ites
# gt and dgt columns are removed; gb (col B) and btes (col N) are newly added.

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Insert new column B for 'gb' (shifts hp..ites one column right)
    $ws.Columns.Item(2).Insert()

    # Remove the old 'gt' and 'dgt' columns, now sitting at F and G
    $ws.Columns.Item(6).Delete()
    $ws.Columns.Item(6).Delete()

    # Insert new column N for 'btes' (right after ttes, before ites)
    $ws.Columns.Item(14).Insert()

    # Header labels for the two new columns
    $ws.Range("B1").Value = "gb"
    $ws.Range("N1").Value = "btes"
}

# --- Sheet 1 data (row 2) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = 39063.99109145206
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 483537.6274462014
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 2897240.114301849
$ws.Range("F2").Value = 94331.34471502228
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 25342.77928792104
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 4897.602148891324
$ws.Range("N2").Value = 23638.06126801545
$ws.Range("O2").Value = 19940.13531829346

# --- Sheet 2 data (row 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = 30846.52922536713
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1495599.874611417
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 70193.79982138964
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 56602.42752520426
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 8321.538881516826
$ws.Range("N2").Value = 51649.16401227913
$ws.Range("O2").Value = 42574.77934331147

# --- Sheet 3 data (row 2) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = 242452.4252219552
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 943335.270081223
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 1425.925979620855
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 39373.98526588717
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 15130.30412399452
$ws.Range("N2").Value = 53308.16490721726
$ws.Range("O2").Value = 30023.09380555204

# --- Sheet 4 data (row 2) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 11578.49752443177
$ws.Range("O2").Value = 0

# --- Sheet 5 data (row 2) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").Value = 76705.58894163162
$ws.Range("B2").Value = 1930.947398408091
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 28147.3462746636
$ws.Range("O2").Value = 8312.661449003012

# --- Sheet 6 data (row 2) ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
